$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... and its sectors at state and territory level?" ->
#    "... and its sectors? And how do they differ by state and territory?"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "sectors at state and territory level?", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "sectors? And how do they differ by state and territory?", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Trim the "Where did we find the data?" paragraph: drop the sentences
#    describing the dataset mechanics and replace the closing clause.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Center for Systems Science and Engineering. The dataset is a time series that cumulatively counts new cases of the virus daily. From this we can calculate infection rate of covid in Australian states starting from January 2020 until now. We found the data source as it seems to be the only place that publicly publishes covid state",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Center for Systems Science and Engineering. We chose to use this data source as it publishes covid state", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append the new paragraphs at the end of the document (after the
#    "Labour force data..." paragraph, before the section break). Each new
#    paragraph inherits formatting from the one it follows, so only the
#    properties that actually change are toggled (keeps the emitted XML
#    free of redundant w:val="0" overrides).
#    Paragraph preceding this block ("Labour force data...") is plain
#    (bold=False, italic=False).
# ---------------------------------------------------------------------------

# --- Paragraph: "Industry data/real estate*" (italic, not bold) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Italic = $true
$p.Range.Font.ItalicBi = $true
$p.Range.InsertBefore("Industry data/real estate*")

# --- Empty paragraph (bold + italic paragraph mark) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Bold = $true
$p.Range.Font.BoldBi = $true

# --- Paragraph: heading "Describe the data exploration and clean-up process (accompanied by Notebook)" (bold, not italic) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Italic = $false
$p.Range.Font.ItalicBi = $false
$p.Range.InsertBefore("Describe the data exploration and clean-up process (accompanied by Notebook)")

# --- Paragraph: CSSE dataset description (normal, not bold/italic) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Bold = $false
$p.Range.Font.BoldBi = $false
$p.Range.InsertBefore("The CSSE dataset is a time series that cumulatively counts new cases of the virus daily. From this we can calculate infection rate of covid in Australian states starting from January 2020 until now.")

# --- Empty paragraph (bold paragraph mark, not italic) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Bold = $true
$p.Range.Font.BoldBi = $true

# --- Paragraph: "Describe the analysis process." (normal, not bold/italic) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Font.Bold = $false
$p.Range.Font.BoldBi = $false
$p.Range.InsertBefore("Describe the analysis process.")

# --- Trailing empty paragraph (normal paragraph mark, not bold/italic) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
